$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values parse as plain numbers (e.g. "200.95"). The source data keeps
# every Price/Volume cell as literal text, so force those specific cells to Text format
# before writing the new value to prevent Excel from auto-converting them to numbers.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '67.815.37'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '3.581.06'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '200.95'
$ws.Range('E5').Value = '  +7.71%  '
$ws.Range('D6').Value = '571.61'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '3.573.99'
$ws.Range('E7').Value = '  +1.25%  '
$ws.Range('D8').Value = '0.614'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('D11').Value = '59.95'
$ws.Range('E11').Value = '  +11.82%  '
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('E13').Value = '  +9.90%  '
$ws.Range('D14').Value = '10.22'
$ws.Range('E14').Value = '  +5.47%  '
$ws.Range('D15').Value = '4.159.47'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '3.581.01'
$ws.Range('E16').Value = '  +1.67%  '
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '19.04'
$ws.Range('E18').Value = '  +5.19%  '
$ws.Range('D19').Value = '67.542.40'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').Value = '12.21'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').Value = '402.83'
$ws.Range('E22').Value = '  +4.27%  '
$ws.Range('D23').Value = '12.72'
$ws.Range('E23').Value = '  +15.84%  '
$ws.Range('D24').Value = '4.20'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').Value = '84.58'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '12.39'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '3.88'
$ws.Range('E28').Value = '  +9.53%  '
$ws.Range('D29').Value = '9.18'
$ws.Range('E29').Value = '  +4.53%  '
$ws.Range('D30').Value = '7.66'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = '31.40'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').Value = '677.11'
$ws.Range('E32').Value = '  +10.02%  '
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').Value = '63.25'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').Value = '41.33'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  +1.76%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '3.22'
$ws.Range('E39').Value = '  +11.41%  '
$ws.Range('D40').Value = '0.0₃0759'
$ws.Range('E40').Value = '  +3.59%  '
$ws.Range('D41').Value = '3.185.35'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('E44').Value = '  +7.09%  '
$ws.Range('E45').Value = '  +14.18%  '
$ws.Range('D46').Value = '2.78'
$ws.Range('E46').Value = '  +18.95%  '
$ws.Range('D47').Value = '0.0409'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').Value = '3.09'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '8.59'
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').Value = '138.86'
$ws.Range('E51').Value = '  +0.92%  '
